# Disaggregation of commodity Copper
# 1. Rename the "Copper ores and concentrates" row label to "Copper" on every
#    year sheet (2000-2100). The label lives in column C, row 4 of each sheet
#    and all sheets share the same underlying string table entry, so once
#    every sheet's C4 is updated the old text is fully replaced.
# 2. A handful of sheets also carry an updated D4 total (same value, last
#    significant digit rounded differently after the disaggregation).

$wb = $excel.ActiveWorkbook

$years = 2000..2100
foreach ($year in $years) {
    $ws = $wb.Worksheets.Item("$year")
    $ws.Range("C4").Value = "Copper"
}

$d4updates = @{
    2026 = 63821.67851678839
    2041 = 384877.2807602866
    2048 = 1319979.23390106
    2054 = 3270520.887466246
    2058 = 3380988.514497868
    2073 = 1438918.7708966
    2074 = 1696090.278107328
    2090 = 2287182.201352461
}

foreach ($year in $d4updates.Keys) {
    $ws = $wb.Worksheets.Item("$year")
    $ws.Range("D4").Value = $d4updates[$year]
}
